$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "61.969.42"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "3.437.17"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue "D5" "412.61"
$ws.Range("E5").Value = "  +0.77%  "
Set-TextValue "D6" "129.97"
$ws.Range("E6").Value = "  +1.40%  "
Set-TextValue "D7" "0.635"
$ws.Range("E7").Value = "  +2.26%  "
$ws.Range("E8").Value = "  +0.09%  "
Set-TextValue "D9" "0.734"
$ws.Range("E9").Value = "  -1.81%  "
Set-TextValue "D10" "0.141"
$ws.Range("E10").Value = "  +1.16%  "
Set-TextValue "D11" "43.61"
$ws.Range("E11").Value = "  +1.95%  "
Set-TextValue "D12" "0.0000220"
$ws.Range("E12").Value = "  +11.12%  "
Set-TextValue "D13" "9.29"
$ws.Range("E13").Value = "  +5.35%  "
$ws.Range("D14").Value = "3.982.21"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  +0.53%  "
Set-TextValue "D16" "21.19"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "3.423.97"
$ws.Range("E17").Value = "  +0.22%  "
Set-TextValue "D18" "12.72"
$ws.Range("E18").Value = "  +0.31%  "
Set-TextValue "D19" "1.10"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").Value = "62.067.20"
$ws.Range("E20").Value = "  +0.15%  "
Set-TextValue "D21" "490.95"
$ws.Range("E21").Value = "  +22.12%  "
Set-TextValue "D22" "92.78"
$ws.Range("E22").Value = "  +3.03%  "
Set-TextValue "D23" "3.31"
$ws.Range("E23").Value = "  +4.03%  "
Set-TextValue "D24" "13.57"
$ws.Range("E24").Value = "  +1.66%  "
Set-TextValue "D25" "3.41"
$ws.Range("E25").Value = "  +5.17%  "
Set-TextValue "D26" "34.86"
$ws.Range("E26").Value = "  +5.60%  "
Set-TextValue "D27" "9.17"
$ws.Range("E27").Value = "  +7.87%  "
Set-TextValue "D28" "4.80"
$ws.Range("E28").Value = "  +0.13%  "
Set-TextValue "D29" "7.72"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D30" "2.72"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D31" "12.13"
$ws.Range("E31").Value = "  +3.00%  "
$ws.Range("E32").Value = "  -1.69%  "
Set-TextValue "D33" "0.168"
$ws.Range("E33").Value = "  -1.78%  "
Set-TextValue "D34" "42.21"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D35" "1.00"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D36" "58.44"
$ws.Range("E36").Value = "  +11.30%  "
Set-TextValue "D37" "0.0498"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  +3.79%  "
$ws.Range("E39").Value = "  -0.07%  "
Set-TextValue "D40" "150.92"
$ws.Range("E40").Value = "  +7.77%  "
Set-TextValue "D41" "2.17"
$ws.Range("E41").Value = "  +10.57%  "
Set-TextValue "D42" "0.137"
$ws.Range("E42").Value = "  +4.91%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D43" "0.325"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "2.97"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D45" "2.70"
$ws.Range("E45").Value = "  +13.16%  "
Set-TextValue "D46" "4.30"
$ws.Range("E46").Value = "  +7.28%  "
Set-TextValue "D47" "2.42"
$ws.Range("E47").Value = "  +26.48%  "
Set-TextValue "D48" "16.66"
$ws.Range("E48").Value = "  -0.32%  "
Set-TextValue "D49" "22.96"
$ws.Range("E49").Value = "  +5.49%  "
Set-TextValue "D50" "118.85"
$ws.Range("E50").Value = "  +24.21%  "
Set-TextValue "D51" "0.147"
$ws.Range("E51").Value = "  +16.49%  "
